$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Neodymium")
$ws.Range("C1").Value = 2030
$ws.Range("C2").Value = [double]"1.625819899935208E-06"
$ws.Range("D2").Value = [double]"0.005885782825730004"
$ws.Range("E2").Value = [double]"0.006492109438477343"
$ws.Range("B3").Value = [double]"2.183980475909259E-12"
$ws.Range("C3").Value = [double]"7.919668242014081E-05"
$ws.Range("D3").Value = [double]"0.005501068836496137"
$ws.Range("E3").Value = [double]"0.005498836806399566"
$ws.Range("B4").Value = [double]"3.409259119931335E-14"
$ws.Range("C4").Value = [double]"7.154027686980761E-05"
$ws.Range("D4").Value = [double]"0.003927863988361095"
$ws.Range("E4").Value = [double]"0.004603584064156631"
$ws.Range("C5").Value = [double]"1.587093575472107E-09"
$ws.Range("D5").Value = [double]"0.0002009765232936477"
$ws.Range("E5").Value = [double]"0.0003658129441505158"

$ws = $wb.Worksheets.Item("Dysprosium")
$ws.Range("C1").Value = 2030

$ws = $wb.Worksheets.Item("Copper")
$ws.Range("C1").Value = 2030
$ws.Range("B2").Value = [double]"3.278497091721097E-06"
$ws.Range("C2").Value = [double]"0.002541871850645688"
$ws.Range("D2").Value = [double]"0.5379000242912599"
$ws.Range("E2").Value = [double]"0.4848590208778999"
$ws.Range("B3").Value = [double]"2.229370101113288E-05"
$ws.Range("C3").Value = [double]"0.009197341277647557"
$ws.Range("D3").Value = [double]"0.3863510865451887"
$ws.Range("E3").Value = [double]"0.3402597778765012"
$ws.Range("B4").Value = [double]"6.612099022439717E-05"
$ws.Range("C4").Value = [double]"0.002447260422877557"
$ws.Range("D4").Value = [double]"0.2766023524654511"
$ws.Range("E4").Value = [double]"0.3027710330745292"
$ws.Range("B5").Value = [double]"2.076994439830034E-05"
$ws.Range("C5").Value = [double]"0.005407093418076853"
$ws.Range("D5").Value = [double]"0.5124365740692156"
$ws.Range("E5").Value = [double]"0.3561425361412697"

$ws = $wb.Worksheets.Item("Raw silicon")
$ws.Range("C1").Value = 2030
$ws.Range("B2").Value = [double]"4.96652837099915E-07"
$ws.Range("C2").Value = [double]"3.462735458016471E-05"
$ws.Range("D2").Value = [double]"0.01385766416475183"
$ws.Range("E2").Value = [double]"0.0139053799804532"
$ws.Range("B3").Value = [double]"5.30035999530297E-07"
$ws.Range("C3").Value = [double]"0.0001164057132748931"
$ws.Range("D3").Value = [double]"0.006708009055267363"
$ws.Range("E3").Value = [double]"0.00648446477114007"
$ws.Range("B4").Value = [double]"3.397047964529607E-06"
$ws.Range("C4").Value = [double]"3.246107553390979E-05"
$ws.Range("D4").Value = [double]"0.006270289992800229"
$ws.Range("E4").Value = [double]"0.007025531601665832"
$ws.Range("B5").Value = [double]"1.823860200208514E-06"
$ws.Range("C5").Value = [double]"4.124207012744105E-05"
$ws.Range("D5").Value = [double]"0.01240009647204811"
$ws.Range("E5").Value = [double]"0.009771097180398069"
